# Add 2022-Q3 data
#
# 1. Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet (i.e. as the second tab, right after the summary sheet).
# 2. Populate it with the fund holdings table for 2022-Q3 (copying the
#    formatting from the 2022-Q2 sheet so headers / index column keep the
#    same look).
# 3. Insert a new row at the top of the "总计" (summary) sheet's data table
#    for the 2022-Q3 totals, pushing the older quarters down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create + position the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Re-fetch the 2022-Q2 sheet by name (defensive - names are stable refs)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy formatting (fonts / borders / alignment) from the 2022-Q2 table so
# the new sheet's header row + index column match the established style.
$q2Sheet.Range("A1:H14").Copy()
$q3Sheet.Range("A1:H14").PasteSpecial(-4122)
# One extra data row (2022-Q3 has 14 funds vs 13 for 2022-Q2) - reuse the
# last row's formatting for it.
$q2Sheet.Range("A14:H14").Copy()
$q3Sheet.Range("A15:H15").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: fill in the 2022-Q3 header + fund rows
# ---------------------------------------------------------------------
$q3Sheet.Cells.Item(1,2).Value = "基金代码"
$q3Sheet.Cells.Item(1,3).Value = "基金名称"
$q3Sheet.Cells.Item(1,4).Value = "基金规模"
$q3Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q3Sheet.Cells.Item(1,6).Value = "仓位占比"
$q3Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3Sheet.Cells.Item(1,8).Value = "仓位排名"

# index | fund code | fund name | scale | position | ratio | held value | rank
$funds = @(
    @(0,  "512980", "广发中证传媒ETF",             "44.76", "99.29", "3.15", "1.4099", 6),
    @(1,  "159869", "华夏中证动漫游戏ETF",           "6.35", "99.31", "6.67", "0.4235", 4),
    @(2,  "516010", "国泰中证动漫游戏ETF",           "3.78", "97.86", "6.42", "0.2427", 5),
    @(3,  "160629", "鹏华中证传媒指数（LOF）A",       "6.41", "94.58", "2.98", "0.1910", 6),
    @(4,  "161030", "富国中证体育产业指数A",          "1.59", "94.00", "4.62", "0.0735", 6),
    @(5,  "516770", "华泰柏瑞中证动漫游戏ETF",        "0.99", "96.39", "6.59", "0.0652", 4),
    @(6,  "159805", "鹏华中证传媒ETF",               "1.71", "98.37", "3.12", "0.0534", 6),
    @(7,  "164818", "工银瑞信中证传媒指数（LOF）A",    "1.65", "93.46", "2.95", "0.0487", 6),
    @(8,  "517500", "国泰中证沪港深动漫游戏ETF",       "0.53", "92.78", "4.63", "0.0245", 6),
    @(9,  "013278", "富国中证体育产业指数C",          "0.42", "94.00", "4.62", "0.0194", 6),
    @(10, "003397", "银华体育文化灵活配置混合",        "0.32", "81.07", "3.43", "0.0110", 10),
    @(11, "010677", "工银瑞信中证传媒指数（LOF）C",    "0.21", "93.46", "2.95", "0.0062", 6),
    @(12, "015675", "鹏华中证传媒指数（LOF）C",        "0.17", "94.58", "2.98", "0.0051", 6),
    @(13, "516190", "华夏中证文娱传媒ETF",            "0.13", "96.01", "2.57", "0.0033", 9)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = $i + 2
    $f = $funds[$i]
    $q3Sheet.Cells.Item($row, 1).Value = $f[0]
    $q3Sheet.Cells.Item($row, 2).Value = "'" + $f[1]
    $q3Sheet.Cells.Item($row, 3).Value = $f[2]
    $q3Sheet.Cells.Item($row, 4).Value = "'" + $f[3]
    $q3Sheet.Cells.Item($row, 5).Value = "'" + $f[4]
    $q3Sheet.Cells.Item($row, 6).Value = "'" + $f[5]
    $q3Sheet.Cells.Item($row, 7).Value = "'" + $f[6]
    $q3Sheet.Cells.Item($row, 8).Value = $f[7]
}

# ---------------------------------------------------------------------
# Step 3: insert the 2022-Q3 row into the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Restore the index-column style (bold / centered / bordered) on the new
# A2 cell by copying it from A3 (pushed-down former A2, same style).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
# Clear the leftover blank-row formatting PasteSpecial left on B2:D2 so
# the data cells have no explicit style, matching the rest of the table.
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 14
$total.Cells.Item(2,4).Value = 2.58

# Renumber the index column (A) for the rows that got pushed down one
# position: they keep their original values, but need to be 1 greater.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(9,1).Value = 7
